# Actualización automática 2025-12-11 17:30:09
$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" (sheet1): per-category December figures ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

$wsGrupo.Range("O3").Value = 1058.37

$wsGrupo.Range("D7").Value = 475.2
$wsGrupo.Range("K7").Value = 558.36
$wsGrupo.Range("M7").Value = 154.22

$wsGrupo.Range("M18").Value = 1358.35

$wsGrupo.Range("I19").Value = 26.1

$wsGrupo.Range("E24").Value = 307.46

$wsGrupo.Range("D29").Value = 88.13
$wsGrupo.Range("E29").Value = 95.48999999999999
$wsGrupo.Range("O29").Value = 1058.37

$wsGrupo.Range("M36").Value = 12890.39
$wsGrupo.Range("O36").Value = 2645.92

$wsGrupo.Range("I53").Value = 102.6

# Row 56 "count of 54 clients per category" summary labels
$wsGrupo.Range("D56").Value = "7 de 54"
$wsGrupo.Range("E56").Value = "3 de 54"
$wsGrupo.Range("I56").Value = "6 de 54"
$wsGrupo.Range("K56").Value = "4 de 54"
$wsGrupo.Range("M56").Value = "14 de 54"
$wsGrupo.Range("O56").Value = "4 de 54"

# --- Sheet "VENTA MENSUAL" (sheet2): December ("diciembre") column F ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$wsMensual.Range("F3").Value = 1295.73
$wsMensual.Range("F7").Value = 1187.78
$wsMensual.Range("F18").Value = 1358.35
$wsMensual.Range("F19").Value = 2999.48
$wsMensual.Range("F24").Value = 4487.71
$wsMensual.Range("F29").Value = 3788.15
$wsMensual.Range("F36").Value = 16661.64
$wsMensual.Range("F55").Value = 1233.7
$wsMensual.Range("F56").Value = 1233.7
$wsMensual.Range("F60").Value = 58114.96
